# Fruta / hortaliza, semanal
# Insert a new weekly data row above the existing row 74 (pushing the
# existing rows 74-85 down to 75-86) and populate it with the new record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(74).Insert()

$ws.Cells.Item(74, 1).Value = 6
$ws.Cells.Item(74, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(74, 3).Value = "Metropolitana"
$ws.Cells.Item(74, 4).Value = 44642
$ws.Cells.Item(74, 5).Value = 13
$ws.Cells.Item(74, 6).Value = "Fruta"
$ws.Cells.Item(74, 7).Value = 100101
$ws.Cells.Item(74, 8).Value = "Berries"
$ws.Cells.Item(74, 9).Value = 100101008
$ws.Cells.Item(74, 10).Value = "Mora"
$ws.Cells.Item(74, 11).Value = "Sin especificar"
$ws.Cells.Item(74, 12).Value = "Primera"
$ws.Cells.Item(74, 13).Value = 250
$ws.Cells.Item(74, 14).Value = 6000
$ws.Cells.Item(74, 15).Value = 6000
$ws.Cells.Item(74, 16).Value = 6000
$ws.Cells.Item(74, 17).Value = "`$/bandeja 2 kilos"
$ws.Cells.Item(74, 18).Value = "Provincia de Linares"
$ws.Cells.Item(74, 19).Value = 3000
$ws.Cells.Item(74, 20).Value = 2
